$d = $word.ActiveDocument

$bullet1 = [char]0x2022 + " Uncovered decades of demographic miscoding in voter files, discovering 500,000+ previously mischaracterized Democratic voters"
$bullet2 = [char]0x2022 + " Developed Python boundary estimation algorithm enabling mapping and analysis at every level of election in the United States"
$bullet3 = [char]0x2022 + " Algorithm reduced mapping costs by 75%, saving campaigns and organizations `$5M+ and enabling smaller nonprofits to conduct redistricting analysis"

$insertText = "`r" + $bullet1 + "`r" + $bullet2 + "`r" + $bullet3

foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Data, Technology and Strategy Consulting*") {
        $p.Range.InsertAfter($insertText)
        break
    }
}
